# Swap the full row contents (columns B through AC) between pairs of rows.
# Column A (the sequential id) stays fixed per physical row; only the
# match data (id/odds/etc.) moves between the two rows of each pair.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(8, 9),
    @(20, 21),
    @(26, 27),
    @(30, 31),
    @(43, 44),
    @(54, 55),
    @(56, 57)
)

$firstCol = 2   # column B
$lastCol  = 29  # column AC

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}
